# Journal de travail - add two new entries (rows 37 & 38) to Tableau1
# and refresh the dependent views (table range, dimension, selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item(1)

# Grow the table by two rows; this keeps the table definition,
# autofilter and worksheet dimension in sync automatically.
$newRow1 = $tbl.ListRows.Add()
$newRow2 = $tbl.ListRows.Add()

# --- Row 37 -----------------------------------------------------------
# 09.05.2022 - Réalisation - 1h15 - Adaptation du template au format MVC
# en PHP - Revu - c90e884f414904d00272d3cbe9c300baa29374c0
$ws.Range("A37").Value = 44690
$ws.Range("B37").Value = "Réalisation"
$ws.Range("C37").Value = 1.25
$ws.Range("D37").Value = "Adaptation du template au format MVC en PHP"
$ws.Range("E37").Value = "Revu"
$ws.Range("F37").Value = "c90e884f414904d00272d3cbe9c300baa29374c0"

# --- Row 38 -----------------------------------------------------------
# 09.05.2022 - Analyse - 0h30 - Continuation de la documentation
$ws.Range("A38").Value = 44690
$ws.Range("B38").Value = "Analyse"
$ws.Range("C38").Value = 0.5
$ws.Range("D38").Value = "Continuation de la documentation"

# Copy the formatting (date number format, wrap text, ...) from the
# previous last row so the new rows look the same as the rest of the
# table. Formats are applied after the values so the recalculation of
# dependent formulas (e.g. the SUM in H7) is not skipped.
$ws.Range("A36:F36").Copy()
$ws.Range("A37:F37").PasteSpecial(-4122)

$ws.Range("A36:D36").Copy()
$ws.Range("A38:D38").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Put the selection where it ends up after typing the last entry.
$ws.Range("F38").Select() | Out-Null
